$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column F (dSF)
$updates = @{
    2  = -2
    3  = 3
    5  = -6
    6  = 3
    7  = 2
    8  = -1
    9  = 3
    10 = -1
    11 = 4
    12 = 4
    13 = 3
    14 = -2
    15 = 1
    16 = -1
    17 = -1
    19 = -6
    20 = 1
    21 = 3
    22 = 1
    23 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
